$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ifoCAST full-series evaluation data: the rolling window shifts up by
# one row (oldest row 2 drops off) and a new row of values is appended at
# row 11.
$data = @{
    2  = @(-0.005035358036557557, 0.3637990656943072, 0.1855351867275108, 0.4307379559865961, 0.4469673578304695, 14)
    3  = @(-0.01898232632975465,  0.3165751929851393, 0.1363897707384557, 0.3693098573534908, 0.3838817902180699, 13)
    4  = @(0.00780077805212256,   0.3794346132818944, 0.1847968220473418, 0.4298800088947401, 0.4489210943938488, 12)
    5  = @(0.03532231998103826,   0.2564852750112934, 0.1004647918617584, 0.3169618145167623, 0.3303616777566779, 11)
    6  = @(0.002000156080238219,  0.3449448151542904, 0.1826303857083173, 0.4273527649475516, 0.4504644332784307, 10)
    7  = @(-0.05875859174690282,  0.3248866822092542, 0.1283271032654252, 0.3582277254281489, 0.3748117494910135, 9)
    8  = @(-0.02790046359007027,  0.3120137538489823, 0.1287774372279015, 0.3588557331684997, 0.3824716629792676, 8)
    9  = @(-0.03188629724616485,  0.2657389849834738, 0.1191175486163901, 0.3451341023665874, 0.3711930514458828, 7)
    10 = @(-0.0881379738596985,   0.4006321887415465, 0.2481252247171507, 0.4981216966938408, 0.5370552175868598, 6)
    11 = @(0.003999580097293598,  0.1842291207164398, 0.04734795102614385, 0.2175958433108129, 0.2432384488140937, 5)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
